# "analisis de metas.xlsx" -- creada la estructura de la base de datos
# Adds a daily-goal ("Meta") input area above the existing table, makes the
# F (Meta diaria) column formula-driven off that input, bumps the working
# month from December to November, and gives the "Outs" banner its own
# larger-font title row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Meta" control cells in row 2 (E2 label, F2 the daily goal input) ---
$ws.Range("E2").Value = "Meta"
$ws.Range("F2").Value = 40

# --- Title banner: move "Outs" from B3 to A3, give it a big font, and size
#     the row to fit it ---
$ws.Range("B3").ClearContents()
$ws.Range("A3").Value = "Outs"
$ws.Range("A3").Font.Size = 28
$ws.Rows.Item(3).RowHeight = 36

# --- Working month changes from December (12) to November (11) ---
$ws.Range("B5").Value = 11

# --- F (Meta diaria) column becomes a formula driven by the new F2 input,
#     instead of hard-coded numbers. F5 is the anchor; F6:F35 fill down. ---
$ws.Range("F5").Formula = "=`$F`$2*C5"
$ws.Range("F6:F35").Formula = "=`$F`$2*C6"

# --- Printing: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Match the author's final selection ---
$ws.Range("F3").Select()
